# Add a new "BAM Sample ID:" header column (column S) to the accessioning
# template, mirroring the style of the existing header row cells.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New header cell S1, with the same bold header styling as the other
# header cells in row 1 (e.g. R1 / "Genome Build:").
$ws.Range("S1").Value = "BAM Sample ID:"
$ws.Range("S1").Font.Bold = $true

# Match the new column's width to the diff (stored width of 15).
# Excel's ColumnWidth property is in "characters" and gets converted to
# the stored spreadsheetML width using the workbook's max-digit-width,
# so 14.17 characters here round-trips to a stored width of 15.
$ws.Range("S1").ColumnWidth = 14.17

# Update the active selection to match the post-edit state.
$ws.Range("S2").Select()
